$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '47.298.70'
Set-TextValue 'E2' '  +0.24%  '
Set-TextValue 'D3' '2.490.04'
Set-TextValue 'E3' '  -0.35%  '
Set-TextValue 'E4' '  +0.14%  '
Set-TextValue 'D5' '321.13'
Set-TextValue 'E5' '  -0.76%  '
Set-TextValue 'D6' '108.41'
Set-TextValue 'E6' '  +2.74%  '
Set-TextValue 'D7' '0.521'
Set-TextValue 'E7' '  +0.03%  '
Set-TextValue 'E8' '  +0.07%  '
Set-TextValue 'D9' '0.533'
Set-TextValue 'E9' '  -1.50%  '
Set-TextValue 'E10' '  +4.19%  '
Set-TextValue 'D11' '0.0809'
Set-TextValue 'E11' '  -0.58%  '
Set-TextValue 'E12' '  +0.23%  '
Set-TextValue 'D13' '18.33'
Set-TextValue 'E13' '  +0.03%  '
Set-TextValue 'E14' '  -1.90%  '
Set-TextValue 'D15' '2.879.84'
Set-TextValue 'E15' '  +0.11%  '
Set-TextValue 'D16' '2.488.83'
Set-TextValue 'E16' '  -2.28%  '
Set-TextValue 'E17' '  +0.07%  '
Set-TextValue 'D18' '47.212.67'
Set-TextValue 'E18' '  +0.44%  '
Set-TextValue 'D19' '12.88'
Set-TextValue 'E19' '  +1.35%  '
Set-TextValue 'D20' '6.60'
Set-TextValue 'E20' '  -0.20%  '
Set-TextValue 'E21' '  -0.44%  '
Set-TextValue 'D22' '2.66'
Set-TextValue 'E22' '  +12.65%  '
Set-TextValue 'D23' '70.29'
Set-TextValue 'E23' '  -0.46%  '
Set-TextValue 'E24' '  -2.28%  '
Set-TextValue 'D25' '2.56'
Set-TextValue 'E25' '  +0.16%  '
Set-TextValue 'E26' '  +0.08%  '
Set-TextValue 'D27' '25.73'
Set-TextValue 'E27' '  -2.01%  '
Set-TextValue 'E28' '  +3.57%  '
Set-TextValue 'E29' '  -0.52%  '
Set-TextValue 'B30' 'Kaspa'
Set-TextValue 'C30' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D30' '0.137'
Set-TextValue 'E30' '  +1.70%  '
Set-TextValue 'B31' 'InjectiveProtocol'
Set-TextValue 'C31' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D31' '34.52'
Set-TextValue 'E31' '  -1.46%  '
Set-TextValue 'D32' '49.64'
Set-TextValue 'E32' '  +0.21%  '
Set-TextValue 'D33' '20.74'
Set-TextValue 'E33' '  +5.09%  '
Set-TextValue 'D34' '5.33'
Set-TextValue 'E34' '  +0.14%  '
Set-TextValue 'E35' '  +0.24%  '
Set-TextValue 'D36' '1.00'
Set-TextValue 'E36' '  +0.32%  '
Set-TextValue 'E37' '  +1.28%  '
Set-TextValue 'D38' '4.67'
Set-TextValue 'E38' '  +1.37%  '
Set-TextValue 'E39' '  -1.52%  '
Set-TextValue 'D40' '23.26'
Set-TextValue 'E40' '  +8.36%  '
Set-TextValue 'E41' '  -0.49%  '
Set-TextValue 'E42' '  +0.03%  '
Set-TextValue 'D43' '117.74'
Set-TextValue 'E43' '  -3.96%  '
Set-TextValue 'E44' '  -0.30%  '
Set-TextValue 'D45' '1.986.81'
Set-TextValue 'D46' '3.04'
Set-TextValue 'E46' '  +1.22%  '
Set-TextValue 'E47' '  -5.97%  '
Set-TextValue 'E48' '  -0.09%  '
Set-TextValue 'E49' '  -1.74%  '
Set-TextValue 'E50' '  -5.96%  '
Set-TextValue 'D51' '56.63'
Set-TextValue 'E51' '  +3.79%  '
